$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the old "note" column
# (B) to column C (carrying its text + the custom column width/bestFit
# formatting with it), and leaves a blank column B for the new "item #"
# data.
$ws.Columns.Item(2).Insert()

# --- Header row ---
$ws.Range("A1").Value = "Checkin meeting"
$ws.Range("B1").Value = "#"
# C1 already holds "CAPSTONE CHECKIN NOTES" after the column shift.

# --- Meeting 2 rows (previously rows 2-9) ---
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Define HOW we are measuring? Sensitivity"
# D2 held the old stray "sensitivity" value that was shifted over from the
# old column C; it is no longer needed now that it has been merged into C2.
$ws.Range("D2").ClearContents()

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 4

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 5

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 6

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 7

$ws.Range("A9").Value = 2
$ws.Range("B9").Value = 8

# --- Meeting 3 rows (new) ---
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = "Should we have a full data dictionary in readme on original data? Or post cleanup/dropping columns?"

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 10

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 11

$ws.Range("A13").Value = 3
$ws.Range("B13").Value = 12

$ws.Range("A14").Value = 3
$ws.Range("B14").Value = 13

$ws.Range("A15").Value = 3
$ws.Range("B15").Value = 14

$ws.Range("A16").Value = 3

# Update the selection to match the committed workbook state.
$ws.Range("C9").Select()
